$wb = $excel.ActiveWorkbook

# --- Update the BAU user id (xx9381 -> xx0770) on every sheet's M2/N2 cells ---
$sheetNames = @("PSPIE", "INTA", "RD", "MPIE", "UPIE")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("M2").Value = "xx0770"
    $ws.Range("N2").Value = "xx0770"
}

# --- Widen column Q on UPIE to fit the longer value (credit risk page handling) ---
$upie = $wb.Worksheets.Item("UPIE")
$upie.Columns("Q").ColumnWidth = 19.5

# --- Update cached selections on each sheet ---
# PSPIE keeps its existing N2 selection (unchanged by the edit).
$pspie = $wb.Worksheets.Item("PSPIE")
$pspie.Range("N2").Select()

$inta = $wb.Worksheets.Item("INTA")
$inta.Range("M2:N2").Select()

$rd = $wb.Worksheets.Item("RD")
$rd.Range("N2").Select()

$mpie = $wb.Worksheets.Item("MPIE")
$mpie.Range("N2").Select()

# UPIE is selected last so it remains the active/visible tab, matching the
# original workbook where UPIE was the tabSelected sheet.
$upie.Range("N2").Select()
